$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price (D) and Volume(1h) (E) columns store plain text (e.g.
# thousands-dot-grouped prices, pre-formatted percentage strings with
# padding spaces). Force the whole data range to Text format first so
# Excel does not auto-coerce numeric-looking values (like "306.50")
# into real numbers, which would silently drop significant trailing
# zeros and change the cell type away from text.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '26.469.67'
$ws.Range("E2").Value = '  -2.03%  '

$ws.Range("D3").Value = '1.788.37'
$ws.Range("E3").Value = '  -0.53%  '

$ws.Range("E4").Value = '  +0.25%  '

$ws.Range("B5").Value = 'USDC'
$ws.Range("C5").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("D5").Value = '1.002'
$ws.Range("E5").Value = '  +0.12%  '

$ws.Range("B6").Value = 'BNB'
$ws.Range("C6").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range("D6").Value = '306.50'
$ws.Range("E6").Value = '  -0.26%  '

$ws.Range("D7").Value = '0.4268'
$ws.Range("E7").Value = '  +1.92%  '

$ws.Range("D8").Value = '0.3614'
$ws.Range("E8").Value = '  +1.19%  '

$ws.Range("D9").Value = '0.07154'
$ws.Range("E9").Value = '  +0.95%  '

$ws.Range("D10").Value = '0.8517'
$ws.Range("E10").Value = '  +0.75%  '

$ws.Range("D11").Value = '20.51'
$ws.Range("E11").Value = '  +1.80%  '

$ws.Range("D12").Value = '1.862.50'
$ws.Range("E12").Value = '  +5.62%  '

$ws.Range("D13").Value = '6.496'
$ws.Range("E13").Value = '  +2.25%  '

$ws.Range("D14").Value = '5.266'
$ws.Range("E14").Value = '  -0.40%  '

$ws.Range("D15").Value = '0.06873'
$ws.Range("E15").Value = '  +1.73%  '

$ws.Range("D16").Value = '1.008'
$ws.Range("E16").Value = '  +0.31%  '

$ws.Range("D17").Value = '79.69'
$ws.Range("E17").Value = '  -0.03%  '

$ws.Range("D18").Value = '0.000008855'
$ws.Range("E18").Value = '  +1.93%  '

$ws.Range("D19").Value = '1.001'
$ws.Range("E19").Value = '  -0.02%  '

$ws.Range("D20").Value = '14.98'
$ws.Range("E20").Value = '  -0.27%  '

$ws.Range("D21").Value = '26.484.57'
$ws.Range("E21").Value = '  -1.55%  '

$ws.Range("D22").Value = '5.137'
$ws.Range("E22").Value = '  +1.75%  '

$ws.Range("D23").Value = '11.01'
$ws.Range("E23").Value = '  +0.66%  '

$ws.Range("D24").Value = '2.048.02'
$ws.Range("E24").Value = '  +2.79%  '

$ws.Range("D25").Value = '152.08'
$ws.Range("E25").Value = '  -0.43%  '

$ws.Range("D26").Value = '1.822'
$ws.Range("E26").Value = '  -5.75%  '

$ws.Range("D27").Value = '18.16'
$ws.Range("E27").Value = '  +0.40%  '

$ws.Range("D28").Value = '5.136'
$ws.Range("E28").Value = '  +2.56%  '

$ws.Range("D29").Value = '1.896'
$ws.Range("E29").Value = '  +15.65%  '

$ws.Range("D30").Value = '114.75'
$ws.Range("E30").Value = '  +1.58%  '

$ws.Range("D31").Value = '0.08882'
$ws.Range("E31").Value = '  -0.81%  '

$ws.Range("D32").Value = '0.7436'
$ws.Range("E32").Value = '  +3.36%  '

$ws.Range("E33").Value = '  +5.96%  '

$ws.Range("D34").Value = '4.341'
$ws.Range("E34").Value = '  +1.05%  '

$ws.Range("E35").Value = '  -3.48%  '

$ws.Range("E36").Value = '  +0.25%  '

$ws.Range("D37").Value = '1.119'
$ws.Range("E37").Value = '  +4.23%  '

$ws.Range("D38").Value = '0.05151'
$ws.Range("E38").Value = '  +0.90%  '

$ws.Range("D39").Value = '0.01897'
$ws.Range("E39").Value = '  -0.13%  '

$ws.Range("D40").Value = '0.4973'
$ws.Range("E40").Value = '  +0.63%  '

$ws.Range("D41").Value = '0.1621'
$ws.Range("E41").Value = '  -0.29%  '

$ws.Range("D42").Value = '2.604'
$ws.Range("E42").Value = '  +1.19%  '

$ws.Range("D43").Value = '6.403'
$ws.Range("E43").Value = '  +7.66%  '

$ws.Range("D44").Value = '8.212'
$ws.Range("E44").Value = '  +2.09%  '

$ws.Range("D45").Value = '105.44'
$ws.Range("E45").Value = '  +0.98%  '

$ws.Range("D46").Value = '10.24'
$ws.Range("E46").Value = '  +1.29%  '

$ws.Range("D48").Value = '1.636'
$ws.Range("E48").Value = '  +2.31%  '

$ws.Range("D49").Value = '0.4516'
$ws.Range("E49").Value = '  -0.20%  '

$ws.Range("D50").Value = '0.06200'
$ws.Range("E50").Value = '  -1.51%  '

$ws.Range("D51").Value = '1.762'
$ws.Range("E51").Value = '  +4.27%  '
